$wb = $excel.ActiveWorkbook

# Insert a brand-new sheet "Ark2" in front of the existing "Ark1" sheet.
# Worksheets.Add() inserts before the active sheet by default, which is
# exactly the position we need (first tab, becomes the active/selected one).
$ark2 = $wb.Worksheets.Add()
$ark2.Name = "Ark2"

# Output header cells for the two advisors.
$ark2.Range("A1").Value = "Vejleder 1:"
$ark2.Range("B1").Value = "Vejleder 2:"
